$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.128.17"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.872.00"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.24"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5137"
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3878"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08384"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.114"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.65"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.193"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.60"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.870.35"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.296"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001105"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.01"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06660"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.032"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "28.148.30"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.09"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.247"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").Value = "2.085.58"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.478"
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.54"
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.58"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.98"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1063"
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.040"
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.897"
$ws.Range("E33").Value = "  +4.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.596"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.538"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02436"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06537"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2185"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.205"
$ws.Range("E39").Value = "  -2.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6491"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.998"
$ws.Range("E41").Value = "  +2.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.227"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.34"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6082"
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.02"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.673"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.276"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.011"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.216"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.40"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.82"
$ws.Range("E51").Value = "  -3.58%  "
